$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.560.56"
$ws.Range("E2").Value = "  -1.24%  "

# Row 3
$ws.Range("D3").Value = "1.588.34"
$ws.Range("E3").Value = "  -2.69%  "

# Row 4
$ws.Range("E4").Value = "  +0.45%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.06"
$ws.Range("E5").Value = "  -2.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -3.50%  "

# Row 7
$ws.Range("E7").Value = "  +0.47%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.21"
$ws.Range("E8").Value = "  -5.67%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.251"
$ws.Range("E9").Value = "  -2.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0590"
$ws.Range("E10").Value = "  -3.65%  "

# Row 11
$ws.Range("E11").Value = "  -2.00%  "

# Row 12
$ws.Range("D12").Value = "1.820.42"
$ws.Range("E12").Value = "  -2.26%  "

# Row 13
$ws.Range("D13").Value = "1.639.65"
$ws.Range("E13").Value = "  +0.42%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.85"
$ws.Range("E14").Value = "  -4.41%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.533"
$ws.Range("E15").Value = "  -5.15%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "27.548.77"
$ws.Range("E16").Value = "  -1.24%  "

# Row 17
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.11"
$ws.Range("E17").Value = "  -3.55%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "216.58"
$ws.Range("E18").Value = "  -5.70%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.34"
$ws.Range("E19").Value = "  -4.47%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0691"
$ws.Range("E20").Value = "  -4.05%  "

# Row 21
$ws.Range("E21").Value = "  +0.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.16"
$ws.Range("E22").Value = "  -3.96%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.55"
$ws.Range("E23").Value = "  -4.99%  "

# Row 24
$ws.Range("E24").Value = "  -3.25%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.87"
$ws.Range("E25").Value = "  -1.22%  "

# Row 26
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.46%  "

# Row 27
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.74"
$ws.Range("E27").Value = "  -2.21%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.05"
$ws.Range("E28").Value = "  -2.97%  "

# Row 29
$ws.Range("E29").Value = "  -4.51%  "

# Row 30
$ws.Range("E30").Value = "  -1.94%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0465"
$ws.Range("E31").Value = "  -3.36%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.25"
$ws.Range("E32").Value = "  -4.84%  "

# Row 33
$ws.Range("D33").Value = "1.371.61"
$ws.Range("E33").Value = "  -1.47%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.93"
$ws.Range("E34").Value = "  -5.84%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.51"
$ws.Range("E35").Value = "  -4.53%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.963"
$ws.Range("E36").Value = "  -5.78%  "

# Row 37
$ws.Range("E37").Value = "  -1.28%  "

# Row 38
$ws.Range("E38").Value = "  -3.81%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.537"
$ws.Range("E39").Value = "  -3.59%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("E40").Value = "  -4.80%  "

# Row 41
$ws.Range("E41").Value = "  +0.43%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.966"
$ws.Range("E42").Value = "  -4.74%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.20"
$ws.Range("E43").Value = "  +2.49%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.33"
$ws.Range("E44").Value = "  -1.86%  "

# Row 45
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.77"
$ws.Range("E45").Value = "  -3.54%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.75"
$ws.Range("E46").Value = "  -2.92%  "

# Row 47
$ws.Range("D47").Value = "1.729.57"
$ws.Range("E47").Value = "  -2.41%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.04"
$ws.Range("E48").Value = "  -1.83%  "

# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0967"
$ws.Range("E49").Value = "  -4.68%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0497"
$ws.Range("E50").Value = "  -1.41%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0964"
$ws.Range("E51").Value = "  -7.11%  "
